# Add two new columns ("Relative Sigma" and "Detector to Foil [cm]") to
# Sheet1, inserted right before the existing "Src [src s^-1]" column (F),
# pushing every later column two positions to the right. Populate the new
# columns with the default values used across the existing rows, then move
# the active selection (matching the author's final cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert two blank columns at F:G - everything from old column F onward
# (Src, Rx, Density, Volume, Half-Life, AW, Isotopic Fraction, Gamma E,
# sigma(rel), Lambda) shifts right by two columns.
$ws.Columns("F:G").Insert()

# New column headers.
$ws.Range("F1").Value = "Relative Sigma"
$ws.Range("G1").Value = "Detector to Foil [cm]"

# New column body values for the 7 data rows (rows 2-8).
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 6).Value = 0.01
    $ws.Cells.Item($r, 7).Value = 1
}

# Match the number formatting/style used by the rest of the data columns
# (copy format only from an existing data cell in the same rows).
$ws.Range("E2").Copy()
$ws.Range("F2:G8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Re-assert the Lambda [s^-1] formula (now column R) as one range write so
# it stays a single shared formula across R2:R8, same as before the insert.
$ws.Range("R2:R8").Formula = "=0.6931471806/Q2"

# Restore the worksheet's active cell/selection.
$ws.Range("P19").Select() | Out-Null
